$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.155562
$ws.Range("H2").Value = 0.466686
$ws.Range("I2").Value = 0.04499457894025669
$ws.Range("J2").Value = 0.05066525794134161
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 2.618408666666667
$ws.Range("N2").Value = 7.855226
$ws.Range("O2").Value = 0.004311757517128344
$ws.Range("P2").Value = 0.004320739451956157
$ws.Range("Q2").Value = 0.407324889004
$ws.Range("R2").Value = 3.665924001036
$ws.Range("S2").Value = 0.0001940057139756765
$ws.Range("T2").Value = 0.0002189113788306897
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.155562
$ws.Range("H3").Value = 0.466686
$ws.Range("I3").Value = 0.04499457894025669
$ws.Range("J3").Value = 0.05066525794134161
$ws.Range("M3").Value = 348.4578143333333
$ws.Range("N3").Value = 1045.373443
$ws.Range("O3").Value = 0.5738086722217269
$ws.Range("P3").Value = 0.5750039880707877
$ws.Range("Q3").Value = 54.206794513322
$ws.Range("R3").Value = 487.861150619898
$ws.Range("S3").Value = 0.02581827959888437
$ws.Range("T3").Value = 0.02913272537290658
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.155562
$ws.Range("H4").Value = 0.466686
$ws.Range("I4").Value = 0.04499457894025669
$ws.Range("J4").Value = 0.05066525794134161
$ws.Range("M4").Value = 3.787182
$ws.Range("N4").Value = 7.574364
$ws.Range("O4").Value = 0.006236387262657937
$ws.Range("P4").Value = 0.004166252296022603
$ws.Range("Q4").Value = 0.589141606284
$ws.Range("R4").Value = 3.534849637704
$ws.Range("S4").Value = 0.0002806036189916739
$ws.Range("T4").Value = 0.0002110842472266919
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.155562
$ws.Range("H5").Value = 0.466686
$ws.Range("I5").Value = 0.04499457894025669
$ws.Range("J5").Value = 0.05066525794134161
$ws.Range("M5").Value = 252.408376
$ws.Range("N5").Value = 757.225128
$ws.Range("O5").Value = 0.4156431829984868
$ws.Range("P5").Value = 0.4165090201812336
$ws.Range("Q5").Value = 39.265151787312
$ws.Range("R5").Value = 353.386366085808
$ws.Range("S5").Value = 0.01870169000840497
$ws.Range("T5").Value = 0.02110253694237766
$ws.Range("I6").Value = 0.5057522141597954
$ws.Range("J6").Value = 0.5694923030357509
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 2.618408666666667
$ws.Range("N6").Value = 7.855226
$ws.Range("O6").Value = 0.004311757517128344
$ws.Range("P6").Value = 0.004320739451956157
$ws.Range("Q6").Value = 4.578450767806889
$ws.Range("R6").Value = 41.206056910262
$ws.Range("S6").Value = 0.002180680911207802
$ws.Range("T6").Value = 0.00246062786131194
$ws.Range("I7").Value = 0.5057522141597954
$ws.Range("J7").Value = 0.5694923030357509
$ws.Range("M7").Value = 348.4578143333333
$ws.Range("N7").Value = 1045.373443
$ws.Range("O7").Value = 0.5738086722217269
$ws.Range("P7").Value = 0.5750039880707877
$ws.Range("Q7").Value = 609.3002088989267
$ws.Range("R7").Value = 5483.701880090341
$ws.Range("S7").Value = 0.2902050064802307
$ws.Range("T7").Value = 0.3274603454211743
$ws.Range("I8").Value = 0.5057522141597954
$ws.Range("J8").Value = 0.5694923030357509
$ws.Range("M8").Value = 3.787182
$ws.Range("N8").Value = 7.574364
$ws.Range("O8").Value = 0.006236387262657937
$ws.Range("P8").Value = 0.004166252296022603
$ws.Range("Q8").Value = 6.622123794678
$ws.Range("R8").Value = 39.732742768068
$ws.Range("S8").Value = 0.003154066666447197
$ws.Range("T8").Value = 0.002372648615089897
$ws.Range("I9").Value = 0.5057522141597954
$ws.Range("J9").Value = 0.5694923030357509
$ws.Range("M9").Value = 252.408376
$ws.Range("N9").Value = 757.225128
$ws.Range("O9").Value = 0.4156431829984868
$ws.Range("P9").Value = 0.4165090201812336
$ws.Range("Q9").Value = 441.3517788914373
$ws.Range("R9").Value = 3972.166010022936
$ws.Range("S9").Value = 0.2102124601019098
$ws.Range("T9").Value = 0.2371986811381748
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.2140146666666667
$ws.Range("H10").Value = 0.6420440000000001
$ws.Range("I10").Value = 0.06190136288879069
$ws.Range("J10").Value = 0.06970280846155817
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 2.618408666666667
$ws.Range("N10").Value = 7.855226
$ws.Range("O10").Value = 0.004311757517128344
$ws.Range("P10").Value = 0.004320739451956157
$ws.Range("Q10").Value = 0.5603778579937778
$ws.Range("R10").Value = 5.043400721944001
$ws.Range("S10").Value = 0.0002669036667562327
$ws.Range("T10").Value = 0.0003011676744319979
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.2140146666666667
$ws.Range("H11").Value = 0.6420440000000001
$ws.Range("I11").Value = 0.06190136288879069
$ws.Range("J11").Value = 0.06970280846155817
$ws.Range("M11").Value = 348.4578143333333
$ws.Range("N11").Value = 1045.373443
$ws.Range("O11").Value = 0.5738086722217269
$ws.Range("P11").Value = 0.5750039880707877
$ws.Range("Q11").Value = 74.57508298194355
$ws.Range("R11").Value = 671.1757468374921
$ws.Range("S11").Value = 0.03551953884793227
$ws.Range("T11").Value = 0.04007939284513019
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.2140146666666667
$ws.Range("H12").Value = 0.6420440000000001
$ws.Range("I12").Value = 0.06190136288879069
$ws.Range("J12").Value = 0.06970280846155817
$ws.Range("M12").Value = 3.787182
$ws.Range("N12").Value = 7.574364
$ws.Range("O12").Value = 0.006236387262657937
$ws.Range("P12").Value = 0.004166252296022603
$ws.Range("Q12").Value = 0.8105124933360001
$ws.Range("R12").Value = 4.863074960016
$ws.Range("S12").Value = 0.000386040871060821
$ws.Range("T12").Value = 0.0002903994857921904
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.2140146666666667
$ws.Range("H13").Value = 0.6420440000000001
$ws.Range("I13").Value = 0.06190136288879069
$ws.Range("J13").Value = 0.06970280846155817
$ws.Range("M13").Value = 252.408376
$ws.Range("N13").Value = 757.225128
$ws.Range("O13").Value = 0.4156431829984868
$ws.Range("P13").Value = 0.4165090201812336
$ws.Range("Q13").Value = 54.01909445351468
$ws.Range("R13").Value = 486.171850081632
$ws.Range("S13").Value = 0.02572887950304137
$ws.Range("T13").Value = 0.0290318484562038
$ws.Range("G14").Value = 1.1608855
$ws.Range("H14").Value = 2.321771
$ws.Range("I14").Value = 0.3357732239901092
$ws.Range("J14").Value = 0.2520605430540592
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 2.618408666666667
$ws.Range("N14").Value = 7.855226
$ws.Range("O14").Value = 0.004311757517128344
$ws.Range("P14").Value = 0.004320739451956157
$ws.Range("Q14").Value = 3.039672654207667
$ws.Range("R14").Value = 18.238035925246
$ws.Range("S14").Value = 0.001447772722589772
$ws.Range("T14").Value = 0.001089087932655167
$ws.Range("G15").Value = 1.1608855
$ws.Range("H15").Value = 2.321771
$ws.Range("I15").Value = 0.3357732239901092
$ws.Range("J15").Value = 0.2520605430540592
$ws.Range("M15").Value = 348.4578143333333
$ws.Range("N15").Value = 1045.373443
$ws.Range("O15").Value = 0.5738086722217269
$ws.Range("P15").Value = 0.5750039880707877
$ws.Range("Q15").Value = 404.5196240212588
$ws.Range("R15").Value = 2427.117744127553
$ws.Range("S15").Value = 0.1926695878253731
$ws.Range("T15").Value = 0.1449358174913725
$ws.Range("G16").Value = 1.1608855
$ws.Range("H16").Value = 2.321771
$ws.Range("I16").Value = 0.3357732239901092
$ws.Range("J16").Value = 0.2520605430540592
$ws.Range("M16").Value = 3.787182
$ws.Range("N16").Value = 7.574364
$ws.Range("O16").Value = 0.006236387262657937
$ws.Range("P16").Value = 0.004166252296022603
$ws.Range("Q16").Value = 4.396484669661
$ws.Range("R16").Value = 17.585938678644
$ws.Range("S16").Value = 0.002094011857233507
$ws.Range("T16").Value = 0.001050147816235678
$ws.Range("G17").Value = 1.1608855
$ws.Range("H17").Value = 2.321771
$ws.Range("I17").Value = 0.3357732239901092
$ws.Range("J17").Value = 0.2520605430540592
$ws.Range("M17").Value = 252.408376
$ws.Range("N17").Value = 757.225128
$ws.Range("O17").Value = 0.4156431829984868
$ws.Range("P17").Value = 0.4165090201812336
$ws.Range("Q17").Value = 293.017223776948
$ws.Range("R17").Value = 1758.103342661688
$ws.Range("S17").Value = 0.1395618515849129
$ws.Range("T17").Value = 0.1049854898137958
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = 0.3333333333333333
$ws.Range("G18").Value = 0.1783253333333333
$ws.Range("H18").Value = 0.534976
$ws.Range("I18").Value = 0.05157862002104791
$ws.Range("J18").Value = 0.05807908750729007
$ws.Range("K18").Value = 2
$ws.Range("L18").Value = 0.6666666666666666
$ws.Range("M18").Value = 2.618408666666667
$ws.Range("N18").Value = 7.855226
$ws.Range("O18").Value = 0.004311757517128344
$ws.Range("P18").Value = 0.004320739451956157
$ws.Range("Q18").Value = 0.4669285982862222
$ws.Range("R18").Value = 4.202357384576
$ws.Range("S18").Value = 0.0002223945025988598
$ws.Range("T18").Value = 0.0002509446047263622
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = 0.3333333333333333
$ws.Range("G19").Value = 0.1783253333333333
$ws.Range("H19").Value = 0.534976
$ws.Range("I19").Value = 0.05157862002104791
$ws.Range("J19").Value = 0.05807908750729007
$ws.Range("M19").Value = 348.4578143333333
$ws.Range("N19").Value = 1045.373443
$ws.Range("O19").Value = 0.5738086722217269
$ws.Range("P19").Value = 0.5750039880707877
$ws.Range("Q19").Value = 62.13885589359644
$ws.Range("R19").Value = 559.249703042368
$ws.Range("S19").Value = 0.02959625946930648
$ws.Range("T19").Value = 0.03339570694020405
$ws.Range("E20").Value = 1
$ws.Range("F20").Value = 0.3333333333333333
$ws.Range("G20").Value = 0.1783253333333333
$ws.Range("H20").Value = 0.534976
$ws.Range("I20").Value = 0.05157862002104791
$ws.Range("J20").Value = 0.05807908750729007
$ws.Range("M20").Value = 3.787182
$ws.Range("N20").Value = 7.574364
$ws.Range("O20").Value = 0.006236387262657937
$ws.Range("P20").Value = 0.004166252296022603
$ws.Range("Q20").Value = 0.675350492544
$ws.Range("R20").Value = 4.052102955264
$ws.Range("S20").Value = 0.0003216642489247369
$ws.Range("T20").Value = 0.0002419721316781449
$ws.Range("E21").Value = 1
$ws.Range("F21").Value = 0.3333333333333333
$ws.Range("G21").Value = 0.1783253333333333
$ws.Range("H21").Value = 0.534976
$ws.Range("I21").Value = 0.5057522141597954
$ws.Range("J21").Value = 0.5694923030357509
$ws.Range("M21").Value = 252.408376
$ws.Range("N21").Value = 757.225128
$ws.Range("O21").Value = 0.4156431829984868
$ws.Range("P21").Value = 0.4165090201812336
$ws.Range("Q21").Value = 45.01080778632534
$ws.Range("R21").Value = 405.0972700769281
$ws.Range("S21").Value = 0.02143830180021784
$ws.Range("T21").Value = 0.02419046383068151
